$d = $word.ActiveDocument

# --- Locate the relevant paragraphs robustly (by text, not fixed index) ---
$paraIsPlacedTo = $null
$paraManyToMany = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($null -eq $paraIsPlacedTo -and $t -eq "Is placed to ") {
        $paraIsPlacedTo = $p
    }
    if ($t -eq "   // For many-to-many:") {
        $paraManyToMany = $p
    }
}

if ($null -eq $paraIsPlacedTo) { throw "Could not find 'Is placed to ' paragraph" }
if ($null -eq $paraManyToMany) { throw "Could not find '// For many-to-many:' paragraph" }

# --- Remove the existing _GoBack bookmark ---
# In the original document it sits right at the end of the
# "// For many-to-many:" paragraph (just before its paragraph mark). The
# edit moves it up into the "Is placed to" paragraph instead.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$start = $paraIsPlacedTo.Range.Start
$prefixLen = "Is placed to".Length

# "Is placed to " -> drop the trailing space, leaving "Is placed to"
$trailingSpace = $d.Range($start + $prefixLen, $start + $prefixLen + 1)
$trailingSpace.Text = ""

# Insert ":" (plus a temporary one-character sentinel, see below) right
# after "Is placed to". Toggling a character property on the new text and
# back off again forces Word to keep it in its own run instead of merging
# it back into the previous "Is placed to" run.
$insertPoint = $d.Range($start + $prefixLen, $start + $prefixLen)
$insertPoint.InsertAfter(":~")
$colonRange = $d.Range($start + $prefixLen, $start + $prefixLen + 1)
$colonRange.Bold = 1
$colonRange.Bold = 0

# Bookmark goes right after the colon, before the trailing space. A
# collapsed range landing exactly on a paragraph-mark boundary is not
# reliable here, so we keep a one-character sentinel ("~") after the
# insertion point while the bookmark is created, then remove it.
$bmPos = $start + $prefixLen + 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# Remove the sentinel character and put the trailing space back.
$sentinel = $d.Range($bmPos, $bmPos + 1)
$sentinel.Text = " "
